$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update top summary fields ---
# "VALOR MORA" total value
$ws.Range("E11").Value = 17333
# "Cant. Trabajadores"
$ws.Range("C13").Value = 5
# "Cant. Periodos"
$ws.Range("F13").Value = 3

# --- Fix up the bottom "closing" row style before we delete a row ---
# Row 21 carries the special "last row" border style used to close off the
# table. Copy that formatting onto row 20 first, so that after we remove
# the old row 21 (below), row 20 keeps looking like the closing row.
$ws.Range("B21:J21").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Remove the old trailing data row (its data now lives in row 20); this
# also shifts the signature block below (previously rows 26-27) up to
# rows 25-26.
$ws.Rows.Item(21).Delete()

# --- Rewrite the worker table (rows 16-20) with the updated dataset ---
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "73139580"
$ws.Range("D16").Value = "AMAURY CAFIEL CAMACHO"
$ws.Range("E16").Value = "2402"
$ws.Range("F16").Value = 8667
$ws.Range("G16").Value = 1300000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "73147587"
$ws.Range("D17").Value = "EBERT CAMACHO PEREZ"
$ws.Range("E17").Value = "2403"
$ws.Range("F17").Value = 1733
$ws.Range("G17").Value = 1300000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "9054384"
$ws.Range("D18").Value = "RAFAEL ENRIQUE CAFIEL CONDES"
$ws.Range("E18").Value = "2403"
$ws.Range("F18").Value = 1733
$ws.Range("G18").Value = 1300000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "73144665"
$ws.Range("D19").Value = "RAFAEL ANTONIO CAFIEL CAMACHO"
$ws.Range("E19").Value = "2403"
$ws.Range("F19").Value = 1733
$ws.Range("G19").Value = 1300000

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1002242986"
$ws.Range("D20").Value = "IVAN ALVAREZ PINO"
$ws.Range("E20").Value = "2408"
$ws.Range("F20").Value = 3467
$ws.Range("G20").Value = 1300000

# --- Update the signature block (now at rows 25-26) ---
$ws.Range("B25").Value = "___________________________________"
$ws.Range("H25").Value = "___________________________________"
$ws.Range("B26").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H26").Value = "FIRMA DEL REPRESENTANTE LEGAL"
